$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "Datos" -> "datos"
$ws.Name = "datos"

# Update/correct the refreshed FX values in column G (Oanda amounts)
$ws.Range("G3").Value = 1.06652
$ws.Range("G4").Value = 5163.49
$ws.Range("G5").Value = 0.14518
$ws.Range("G6").Value = 0

# These amounts came back as NaN from the source feed -> blank them out
# (kept as empty-text cells, matching the rest of the blank cells on the sheet)
$ws.Range("C3").Formula = '=""'
$ws.Range("C4").Formula = '=""'
$ws.Range("C5").Formula = '=""'
$ws.Range("C6").Formula = '=""'
$ws.Range("C7").Formula = '=""'
$ws.Range("C8").Formula = '=""'
$ws.Range("C9").Formula = '=""'
$ws.Range("C10").Formula = '=""'
$ws.Range("C11").Formula = '=""'

$ws.Range("G7").Formula = '=""'
$ws.Range("G8").Formula = '=""'
$ws.Range("G9").Formula = '=""'
$ws.Range("G10").Formula = '=""'
$ws.Range("G11").Formula = '=""'

# Drop the stale trailing rows (12-20) that no longer ship with the feed
$ws.Rows("12:20").Delete()
